$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells are formatted as Text so numeric-looking strings
# (e.g. "0.35", "1470569") are preserved as text, matching the source data.
$ws.Range("A7:L8").NumberFormat = "@"
$ws.Range("A12:L14").NumberFormat = "@"
$ws.Range("A17:L20").NumberFormat = "@"
$ws.Range("A22:L30").NumberFormat = "@"

# Row 7
$ws.Range("A7").Value = "Less than Primary"
$ws.Range("B7").Value = "0.35"
$ws.Range("C7").Value = "0.14"
$ws.Range("D7").Value = "0.23"
$ws.Range("E7").Value = "0.09"
$ws.Range("F7").Value = "0.21"
$ws.Range("G7").Value = "0.07"
$ws.Range("H7").Value = "0.07"
$ws.Range("I7").Value = "0.06"
$ws.Range("J7").Value = "0.03"
$ws.Range("K7").Value = "0.01"
$ws.Range("L7").Value = "0.03"

# Row 8
$ws.Range("A8").Value = "Education Completed"
$ws.Range("B8:L8").ClearContents()

# Row 12
$ws.Range("A12").Value = "Household Size"
$ws.Range("B12").Value = "3.35"
$ws.Range("C12").Value = "2.4"
$ws.Range("D12").Value = "3.01"
$ws.Range("E12").Value = "2.49"
$ws.Range("F12").Value = "3.2"
$ws.Range("G12").Value = "2.8"
$ws.Range("H12").Value = "2.74"
$ws.Range("I12").Value = "2.37"
$ws.Range("J12").Value = "2.07"
$ws.Range("K12").Value = "2.02"
$ws.Range("L12").Value = "2.28"

# Row 13
$ws.Range("A13").Value = "Lives Alone"
$ws.Range("B13").Value = "0.15"
$ws.Range("C13").Value = "0.26"
$ws.Range("D13").Value = "0.2"
$ws.Range("E13").Value = "0.23"
$ws.Range("F13").Value = "0.17"
$ws.Range("G13").Value = "0.17"
$ws.Range("H13").Value = "0.16"
$ws.Range("I13").Value = "0.27"
$ws.Range("J13").Value = "0.36"
$ws.Range("K13").Value = "0.23"
$ws.Range("L13").Value = "0.28"

# Row 14
$ws.Range("A14").Value = "Household"
$ws.Range("B14:L14").ClearContents()

# Row 17
$ws.Range("A17").Value = "Less than 15"
$ws.Range("B17").Value = "0.12"
$ws.Range("C17").Value = "0.3"
$ws.Range("D17").Value = "0.07"
$ws.Range("E17").Value = "0.19"
$ws.Range("F17").Value = "0.07"
$ws.Range("G17").Value = "0.08"
$ws.Range("H17").Value = "0.14"
$ws.Range("I17").Value = "-"
$ws.Range("J17").Value = "-"
$ws.Range("K17").Value = "-"
$ws.Range("L17").Value = "-"

# Row 18
$ws.Range("A18").Value = "15 - 24"
$ws.Range("B18").Value = "0.34"
$ws.Range("C18").Value = "0.31"
$ws.Range("D18").Value = "0.19"
$ws.Range("E18").Value = "0.15"
$ws.Range("F18").Value = "0.23"
$ws.Range("G18").Value = "0.21"
$ws.Range("H18").Value = "0.19"
$ws.Range("I18").Value = "-"
$ws.Range("J18").Value = "-"
$ws.Range("K18").Value = "-"
$ws.Range("L18").Value = "-"

# Row 19
$ws.Range("A19").Value = "25 - 49"
$ws.Range("B19").Value = "0.44"
$ws.Range("C19").Value = "0.24"
$ws.Range("D19").Value = "0.52"
$ws.Range("E19").Value = "0.44"
$ws.Range("F19").Value = "0.58"
$ws.Range("G19").Value = "0.51"
$ws.Range("H19").Value = "0.49"
$ws.Range("I19").Value = "-"
$ws.Range("J19").Value = "-"
$ws.Range("K19").Value = "-"
$ws.Range("L19").Value = "-"

# Row 20
$ws.Range("A20").Value = "Age Migrated"
$ws.Range("B20:L20").ClearContents()

# Row 22
$ws.Range("A22").Value = "Before 1965"
$ws.Range("B22").Value = "0.15"
$ws.Range("C22").Value = "0.42"
$ws.Range("D22").Value = "0.07"
$ws.Range("E22").Value = "0.26"
$ws.Range("F22").Value = "0.09"
$ws.Range("G22").Value = "0.12"
$ws.Range("H22").Value = "0.2"
$ws.Range("I22").Value = "-"
$ws.Range("J22").Value = "-"
$ws.Range("K22").Value = "-"
$ws.Range("L22").Value = "-"

# Row 23
$ws.Range("A23").Value = "1965 - 1979"
$ws.Range("B23").Value = "0.48"
$ws.Range("C23").Value = "0.3"
$ws.Range("D23").Value = "0.34"
$ws.Range("E23").Value = "0.35"
$ws.Range("F23").Value = "0.37"
$ws.Range("G23").Value = "0.35"
$ws.Range("H23").Value = "0.32"
$ws.Range("I23").Value = "-"
$ws.Range("J23").Value = "-"
$ws.Range("K23").Value = "-"
$ws.Range("L23").Value = "-"

# Row 24
$ws.Range("A24").Value = "1980 - 1999"
$ws.Range("B24").Value = "0.32"
$ws.Range("C24").Value = "0.18"
$ws.Range("D24").Value = "0.44"
$ws.Range("E24").Value = "0.33"
$ws.Range("F24").Value = "0.49"
$ws.Range("G24").Value = "0.35"
$ws.Range("H24").Value = "0.36"
$ws.Range("I24").Value = "-"
$ws.Range("J24").Value = "-"
$ws.Range("K24").Value = "-"
$ws.Range("L24").Value = "-"

# Row 25
$ws.Range("A25").Value = "After 1999"
$ws.Range("B25").Value = "0.1"
$ws.Range("C25").Value = "0.13"
$ws.Range("D25").Value = "0.19"
$ws.Range("E25").Value = "0.19"
$ws.Range("F25").Value = "0.12"
$ws.Range("G25").Value = "0.22"
$ws.Range("H25").Value = "0.16"
$ws.Range("I25").Value = "-"
$ws.Range("J25").Value = "-"
$ws.Range("K25").Value = "-"
$ws.Range("L25").Value = "-"

# Row 26
$ws.Range("A26").Value = "Migration Cohort"
$ws.Range("B26:L26").ClearContents()

# Row 27
$ws.Range("A27").Value = "Citizen"
$ws.Range("B27").Value = "0.53"
$ws.Range("C27").Value = "-"
$ws.Range("D27").Value = "0.66"
$ws.Range("E27").Value = "0.77"
$ws.Range("F27").Value = "0.63"
$ws.Range("G27").Value = "0.71"
$ws.Range("H27").Value = "0.74"
$ws.Range("I27").Value = "-"
$ws.Range("J27").Value = "-"
$ws.Range("K27").Value = "-"
$ws.Range("L27").Value = "-"

# Row 28
$ws.Range("A28").Value = "English Speakers"
$ws.Range("B28").Value = "0.78"
$ws.Range("C28").Value = "0.94"
$ws.Range("D28").Value = "0.74"
$ws.Range("E28").Value = "0.8"
$ws.Range("F28").Value = "0.86"
$ws.Range("G28").Value = "0.9"
$ws.Range("H28").Value = "0.94"
$ws.Range("I28").Value = "0.99"
$ws.Range("J28").Value = "1"
$ws.Range("K28").Value = "1"
$ws.Range("L28").Value = "1"

# Row 29
$ws.Range("A29").Value = "N"
$ws.Range("B29").Value = "37469"
$ws.Range("C29").Value = "9723"
$ws.Range("D29").Value = "3590"
$ws.Range("E29").Value = "9390"
$ws.Range("F29").Value = "7972"
$ws.Range("G29").Value = "10530"
$ws.Range("H29").Value = "155138"
$ws.Range("I29").Value = "54362"
$ws.Range("J29").Value = "130751"
$ws.Range("K29").Value = "1470569"
$ws.Range("L29").Value = "42940"

# Row 30
$ws.Range("A30").Value = "Acculturation"
$ws.Range("B30:L30").ClearContents()
